$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

# --- Fix up existing "resolutions" descriptions (drop redundant "Number, " prefix) ---
$ws.Range("D22").Value = "Height of screen."
$ws.Range("D23").Value = "Inner width of browser window."
$ws.Range("D24").Value = "Inner width of browser window."

# --- New trackers for user context (CMS-15982) ---
$ws.Range("A26").Value = "user_lastLoginDate"
$ws.Range("B26").Value = "Date"
$ws.Range("C26").Value = "Last time the user logged in."

$ws.Range("A28").Value = "user_creationDate"
$ws.Range("B28").Value = "Date"
$ws.Range("C28").Value = "When the user was created."

$ws.Range("A30").Value = "user_roles"
$ws.Range("B30").Value = "String"
$ws.Range("C30").Value = 'User roles, e.g. "admin, webadmin, editor"'

# Highlight the new tracker names with the same light accent5 fill Excel would
# compute for theme="8" tint="0.79998168889431442" (accent5 5B9BD5 tinted ~0.8)
$ws.Range("A26").Interior.Color = 16247774
$ws.Range("A28").Interior.Color = 16247774
$ws.Range("A30").Interior.Color = 16247774
